$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Environments_OnGoing")

# Column A labels first (matches shared-string insertion order of the source commit)
$ws.Range("A9").Value = "DirectSalesDiscountConsoleByProductBasket"
$ws.Range("A10").Value = "DirectSalesOrderEnrichmentByProductBasket"
$ws.Range("A11").Value = "DirectSalesOrderEnrichmentByOrder"

# Column B URLs, in the order the strings were originally added
$ws.Range("B11").Value = "https://proximus--prxitt--csoe.visualforce.com/apex/apex/NonCommercialSpecifications?orderId="
$ws.Range("B9").Value = "https://proximus--prxitt.lightning.force.com/apex/csdiscounts__DiscountPage?basketId="
$ws.Range("B10").Value = "https://proximus--prxitt--csoe.visualforce.com/apex/apex/NonCommercialSpecifications?basketId="

$ws.Hyperlinks.Add($ws.Range("B9"), $ws.Range("B9").Value)
$ws.Hyperlinks.Add($ws.Range("B11"), $ws.Range("B11").Value)
$ws.Hyperlinks.Add($ws.Range("B10"), $ws.Range("B10").Value)

$ws.Range("B9").Style = "Hyperlink"
$ws.Range("B10").Style = "Hyperlink"
$ws.Range("B11").Style = "Hyperlink"

$ws.Columns.Item(1).AutoFit() | Out-Null

$ws.Range("B10").Select() | Out-Null
